$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12: average of the k column (J)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# New summary labels + formulas in rows 14-17
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Build the bold 12pt / vertically centred look on an unused scratch cell so
# a single compact style record is created, then stamp it onto B14:B17 via
# copy/paste-special (avoids leaving behind extra intermediate style rows).
$tmp = $ws.Range("AA100")
$tmp.Font.Bold = $true
$tmp.Font.Size = 12
$tmp.VerticalAlignment = -4108
$tmp.Copy()
$summaryRange = $ws.Range("B14:B17")
$summaryRange.PasteSpecial(-4122)
$tmp.Clear()

$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

$ws.Range("A14:B17").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
